$wb = $excel.ActiveWorkbook

# Overview sheet: zh-cn (col B) and de-de (col C) status columns
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "In Translation"
$wsOverview.Range("C3").Value = "In Translation"
$wsOverview.Range("B4").Value = "In Translation"
$wsOverview.Range("C4").Value = "In Translation"

# zh-cn sheet: Status column (C)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "In Translation"
$wsZh.Range("C4").Value = "In Translation"

# de-de sheet: Status column (C)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "In Translation"
$wsDe.Range("C4").Value = "In Translation"
